# Rename the original (only) sheet from "Sheet1" to "nbaTeams", and add a
# new blank sheet named "Sheet2" right after it. Also move the selection on
# the nbaTeams sheet from D12 to D21 (matching the target selection), and
# keep nbaTeams as the active/selected tab.

$wb = $excel.ActiveWorkbook

# Rename the existing sheet.
$wb.Sheets(1).Name = "nbaTeams"

# Add a brand-new empty worksheet right after "nbaTeams" and name it "Sheet2".
$newSheet = $wb.Worksheets.Add($null, $wb.Sheets("nbaTeams"))
$newSheet.Name = "Sheet2"

# Make sure "nbaTeams" stays the active tab and update its selection to D21.
[void]$wb.Sheets("nbaTeams").Select()
[void]$wb.Sheets("nbaTeams").Range("D21").Select()
